$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Run 50" results column is inserted where the "Mean" column used to be
# (column AZ), pushing the "Mean" column out to the new last column (BA) and
# recalculating it to include the new Run 50 value.

# Give the new BA1 header cell the same look (bold/border/centered) as the
# other header cells before putting the "Mean" text into it.
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)

$ws.Range("BA1").Value = "Mean"
$ws.Range("AZ1").Value = "Run 50"

# Rows 2-14: column AZ becomes the new Run 50 values, column BA becomes the
# recalculated Mean (same value repeated down the table, as with the other
# per-run columns).
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = 607.00321812
    $ws.Cells.Item($r, 53).Value = 580.21353715
}
